# Update "想去人数" (number of people wanting to go) counts in column F
# across the four worksheets, per the diff between the previous and the
# newly generated gh-pages output.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 51
$ws1.Range("F6").Value = 367
$ws1.Range("F7").Value = 1189
$ws1.Range("F8").Value = 447
$ws1.Range("F9").Value = 7141
$ws1.Range("F10").Value = 84
$ws1.Range("F12").Value = 2052
$ws1.Range("F13").Value = 8007
$ws1.Range("F15").Value = 53
$ws1.Range("F16").Value = 5518
$ws1.Range("F18").Value = 2434
$ws1.Range("F20").Value = 4564
$ws1.Range("F24").Value = 10
$ws1.Range("F25").Value = 387
$ws1.Range("F28").Value = 2394
$ws1.Range("F30").Value = 268
$ws1.Range("F31").Value = 86
$ws1.Range("F32").Value = 157
$ws1.Range("F33").Value = 589
$ws1.Range("F34").Value = 8
$ws1.Range("F36").Value = 1517
$ws1.Range("F38").Value = 6
$ws1.Range("F39").Value = 2351

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 96

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 1281

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 51
$ws4.Range("F4").Value = 1281
$ws4.Range("F7").Value = 96
$ws4.Range("F8").Value = 367
$ws4.Range("F9").Value = 1189
$ws4.Range("F10").Value = 447
$ws4.Range("F11").Value = 7141
$ws4.Range("F12").Value = 84
$ws4.Range("F14").Value = 2052
$ws4.Range("F15").Value = 8007
$ws4.Range("F17").Value = 53
$ws4.Range("F18").Value = 5518
$ws4.Range("F20").Value = 2434
$ws4.Range("F22").Value = 4564
$ws4.Range("F26").Value = 10
$ws4.Range("F28").Value = 387
$ws4.Range("F30").Value = 2394
$ws4.Range("F32").Value = 268
$ws4.Range("F33").Value = 86
$ws4.Range("F34").Value = 157
$ws4.Range("F36").Value = 589
$ws4.Range("F37").Value = 8
$ws4.Range("F40").Value = 1517
$ws4.Range("F42").Value = 6
$ws4.Range("F43").Value = 2351

$wb.Save()
